$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New columns T (plus p rate) and U (ctl rate): headers + per-row formulas
# ---------------------------------------------------------------------------
$ws.Range("T1").Value = "plus p rate"
$ws.Range("U1").Value = "ctl rate"

$ws.Range("T2").Formula = "=Q2-S2/(E2/24)"
$ws.Range("U2").Formula = "=R2-S2/(E2/24)"
$ws.Range("T3").Formula = "=Q3-S3/(E3/24)"
$ws.Range("U3").Formula = "=R3-S3/(E3/24)"

# ---------------------------------------------------------------------------
# 2. Corrected "inc ph" reading for the 120m 2018 p1 trap
# ---------------------------------------------------------------------------
$ws.Range("Q2").Value = 2.56

# ---------------------------------------------------------------------------
# 3. Expand / clarify the footnote about the (still dubious) post-p alk
#    values, and push the two footnote lines far down the sheet (row 5/6 ->
#    row 21/22), matching the author's re-arrangement of the notes block.
# ---------------------------------------------------------------------------
$note1 = $ws.Range("A5").Value2
$note2 = $ws.Range("A6").Value2
$ws.Rows("5:6").Delete()

$ws.Range("A21").Value = "post p alk for 2018 p1 120 180 needs to be checked, right now my old values and I’m not sure where they’re from"
$ws.Range("A22").Value = $note2

# ---------------------------------------------------------------------------
# 4. Light formatting touch-up: re-apply the Normal style on every
#    previously-plain cell (the distinctly-fonted "Courier New" cells are
#    left untouched so they keep standing out).
# ---------------------------------------------------------------------------
$plainCells = @( `
    "A1","B1","C1","D1","E1","F1","G1","H1","I1","J1","K1","L1","M1","N1","O1","P1","Q1","R1","S1", `
    "A2","B2","C2","D2","E2","F2","I2","K2","L2","N2","O2","P2","Q2","R2","S2", `
    "A3","B3","C3","D3","E3","F3","G3","H3","I3","N3","R3","S3", `
    "A21","A22" `
)
foreach ($addr in $plainCells) {
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------------
# 5. Restore the selection / scroll anchor the author left the sheet at.
# ---------------------------------------------------------------------------
[void]$ws.Range("Q3").Select()
